$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.685.31'
$ws.Range("E2").Value = '  -0.07%  '
$ws.Range("D3").Value = '1.900.49'
$ws.Range("E3").Value = '  +0.49%  '
$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = '  -0.54%  '
$ws.Range("D5").Value = "'311.75"
$ws.Range("E5").Value = '  -0.84%  '
$ws.Range("D6").Value = "'0.9992"
$ws.Range("E6").Value = '  -0.56%  '
$ws.Range("D7").Value = "'0.5185"
$ws.Range("E7").Value = '  +7.88%  '
$ws.Range("D8").Value = "'0.3780"
$ws.Range("E8").Value = '  -0.26%  '
$ws.Range("D9").Value = "'0.07236"
$ws.Range("E9").Value = '  -1.35%  '
$ws.Range("D10").Value = "'21.05"
$ws.Range("E10").Value = '  +3.20%  '
$ws.Range("D11").Value = "'0.8941"
$ws.Range("E11").Value = '  -2.65%  '
$ws.Range("D12").Value = "'0.07648"
$ws.Range("E12").Value = '  -0.69%  '
$ws.Range("D13").Value = '1.909.21'
$ws.Range("E13").Value = '  +0.93%  '
$ws.Range("D14").Value = "'5.441"
$ws.Range("E14").Value = '  -0.42%  '
$ws.Range("D15").Value = "'92.04"
$ws.Range("E15").Value = '  +1.19%  '
$ws.Range("D16").Value = "'0.9997"
$ws.Range("E16").Value = '  -0.60%  '
$ws.Range("D17").Value = "'0.000008711"
$ws.Range("E17").Value = '  -0.80%  '
$ws.Range("D18").Value = "'0.9994"
$ws.Range("E18").Value = '  -0.52%  '
$ws.Range("D19").Value = '27.719.06'
$ws.Range("E19").Value = '  -0.18%  '
$ws.Range("D20").Value = "'14.47"
$ws.Range("E20").Value = '  -0.25%  '
$ws.Range("D21").Value = "'5.136"
$ws.Range("E21").Value = '  +0.22%  '
$ws.Range("D22").Value = '2.159.06'
$ws.Range("E22").Value = '  +0.34%  '
$ws.Range("D23").Value = "'10.82"
$ws.Range("E23").Value = '  +0.17%  '
$ws.Range("D24").Value = "'6.578"
$ws.Range("E24").Value = '  -0.08%  '
$ws.Range("D25").Value = "'153.70"
$ws.Range("E25").Value = '  -0.38%  '
$ws.Range("D26").Value = "'1.862"
$ws.Range("E26").Value = '  -2.01%  '
$ws.Range("D27").Value = "'2.183"
$ws.Range("E27").Value = '  +3.26%  '
$ws.Range("D29").Value = "'114.79"
$ws.Range("E29").Value = '  -1.37%  '
$ws.Range("D30").Value = "'4.846"
$ws.Range("E30").Value = '  -1.45%  '
$ws.Range("D31").Value = "'0.08971"
$ws.Range("E31").Value = '  +0.38%  '
$ws.Range("D32").Value = "'3.182"
$ws.Range("E32").Value = '  +0.74%  '
$ws.Range("D33").Value = "'1.237"
$ws.Range("E33").Value = '  +0.16%  '
$ws.Range("D34").Value = "'4.802"
$ws.Range("E34").Value = '  +3.71%  '
$ws.Range("D35").Value = "'0.7774"
$ws.Range("E35").Value = '  +2.29%  '
$ws.Range("D36").Value = "'2.617"
$ws.Range("E36").Value = '  +3.75%  '
$ws.Range("D37").Value = "'0.02087"
$ws.Range("E37").Value = '  +2.64%  '
$ws.Range("D38").Value = "'3.053"
$ws.Range("E38").Value = '  +2.61%  '
$ws.Range("D39").Value = "'1.089"
$ws.Range("E39").Value = '  -0.24%  '
$ws.Range("D40").Value = "'0.5491"
$ws.Range("E40").Value = '  +1.08%  '
$ws.Range("D41").Value = "'0.05255"
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").Value = "'6.668"
$ws.Range("E42").Value = '  -4.03%  '
$ws.Range("D43").Value = "'113.20"
$ws.Range("E43").Value = '  +3.61%  '
$ws.Range("D44").Value = "'8.481"
$ws.Range("E44").Value = '  +2.27%  '
$ws.Range("D45").Value = "'0.1504"
$ws.Range("E45").Value = '  -0.82%  '
$ws.Range("D46").Value = "'0.4784"
$ws.Range("E46").Value = '  +0.22%  '
$ws.Range("D47").Value = "'10.45"
$ws.Range("E47").Value = '  -1.22%  '
$ws.Range("D48").Value = "'0.9994"
$ws.Range("E48").Value = '  -0.47%  '
$ws.Range("D49").Value = "'1.613"
$ws.Range("E49").Value = '  -1.32%  '
$ws.Range("D50").Value = "'66.50"
$ws.Range("E50").Value = '  -1.35%  '
$ws.Range("D51").Value = "'0.05996"
$ws.Range("E51").Value = '  -1.09%  '
